# Applies the CycleTime_s correction (56/40/50/76/80/36/82 -> 300) across the
# Raw_Annotations and Aggregates sheets, and updates which sheet/range is the
# active selection when the workbook is saved.

$wb = $excel.ActiveWorkbook

$wsRaw = $wb.Worksheets.Item("Raw_Annotations")
$wsAgg = $wb.Worksheets.Item("Aggregates")

# --- Raw_Annotations: set CycleTime_s (column E, rows 2:31) to 300 ---
$wsRaw.Range("E2:E31").Value = 300

# --- Aggregates: set CycleTime_s (column C, rows 2:6) to 300 ---
$wsAgg.Range("C2:C6").Value = 300

# Formula cells (F/G/J/K on Raw_Annotations and F/H on Aggregates) depend on
# the edited columns, so Excel recalculates their cached values automatically.
$excel.CalculateFullRebuild()

# --- Update active sheet / selection state ---
# Previously "Aggregates" was the active (tabSelected) sheet with A6 selected.
# Now "Raw_Annotations" becomes the active sheet with E2:E31 selected, and
# Aggregates reverts to a plain (non-selected-tab) view with C2:C6 selected.
$wsAgg.Range("C2:C6").Select()
$wsRaw.Activate()
$wsRaw.Range("E2:E31").Select()

$wb.Save()
